$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# orderReferenceID for row 3 (was "Clone_DIR_C02")
$ws.Range("C3").Value = "DIR_C02"

# OrderId for row 2 (was "51490130")
$ws.Range("Z2").Value = "51490999"
